# Append 5 new rows (302-306) of COVID totals data to the "covid_totals" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append: date, areaType, areaCode, areaName, cumCases, newCases, newDeaths, cumDeaths
$rows = @(
    @("2021-06-09", "overview", "K02000001", "United Kingdom", 4535754, 7540, 6, 127860),
    @("2021-06-10", "overview", "K02000001", "United Kingdom", 4542986, 7393, 7, 127867),
    @("2021-06-11", "overview", "K02000001", "United Kingdom", 4550944, 8125, 17, 127884),
    @("2021-06-12", "overview", "K02000001", "United Kingdom", 4558494, 7738, 12, 127896),
    @("2021-06-13", "overview", "K02000001", "United Kingdom", 4565813, 7490, 8, 127904)
)

$startRow = 302

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds a date formatted as plain text (e.g. "2021-06-09"), not a
    # real Excel date serial, so force the cell to text format before writing.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $data[0]

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
